# Fix Training Data Issue (#48)
#
# The "Date" column (BF) for every data row (rows 2-31) held the value
# "6-24-2007-08", a mangled combination of the game date and the season
# label. Because NBA.com showed the game date one day off from how it
# should be recorded, the data was off by a day. Correct this by writing
# the real, unambiguous ISO-style game date "2008-06-24" into BF2:BF31.
#
# NumberFormat is forced to Text ("@") before writing the value so Excel
# does not auto-detect "2008-06-24" as a date and convert it into a date
# serial number. Afterwards the style is reset back to "Normal" so the
# cells keep using the workbook's default (unstyled) formatting, exactly
# like the original cells, which carried no explicit style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$correctDate = "2008-06-24"
$rng = $ws.Range("BF2:BF31")

$rng.NumberFormat = "@"
$rng.Value = $correctDate
$rng.Style = "Normal"
